# Generate Report for Handback
# Update handback/generation timestamps (and one status code) that were
# refreshed when the handback report was regenerated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
# "Latest HO Xliff Generate Date" for 0505a389-...md was refreshed
$wsOverview.Range("G2").Value = "2016-08-21 14:15:56"
$wsOverview.Range("G4").Value = "2016-08-21 14:15:56"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority changed from "ht" to "mt"
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
# Correspond Handoff Datetime refreshed
$wsZhCn.Range("H2").Value = "2016-08-21 14:15:52"
$wsZhCn.Range("H4").Value = "2016-08-21 14:15:52"
# Correspond Handback DateTime refreshed
$wsZhCn.Range("K2").Value = "2016-08-21 14:16:14"
$wsZhCn.Range("K4").Value = "2016-08-21 14:16:14"

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
# Priority changed from "ht" to "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
# Correspond Handoff Datetime refreshed (shares text with Overview!G2/G4)
$wsDeDe.Range("H2").Value = "2016-08-21 14:15:56"
$wsDeDe.Range("H4").Value = "2016-08-21 14:15:56"
# Correspond Handback DateTime refreshed
$wsDeDe.Range("K2").Value = "2016-08-21 14:16:20"
$wsDeDe.Range("K4").Value = "2016-08-21 14:16:20"
